$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($row, $text)
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-VolumeText {
    param($row, $text)
    $ws.Cells.Item($row, 5).Value = $text
}

Set-PriceText 2 "28.557.76"
Set-VolumeText 2 "  +2.38%  "

Set-PriceText 3 "1.911.25"
Set-VolumeText 3 "  +5.50%  "

Set-VolumeText 4 "  -0.39%  "

Set-PriceText 5 "314.48"
Set-VolumeText 5 "  +1.19%  "

Set-PriceText 6 "1.001"

Set-PriceText 7 "0.5057"
Set-VolumeText 7 "  +1.85%  "

Set-PriceText 8 "0.3964"
Set-VolumeText 8 "  +0.77%  "

Set-PriceText 9 "0.09780"
Set-VolumeText 9 "  +2.63%  "

Set-PriceText 10 "1.161"
Set-VolumeText 10 "  +5.58%  "

Set-PriceText 11 "41.87"

Set-PriceText 12 "6.544"
Set-VolumeText 12 "  +2.01%  "

Set-VolumeText 13 "  +3.52%  "

Set-PriceText 14 "1.906.05"
Set-VolumeText 14 "  +4.91%  "

Set-PriceText 15 "7.582"
Set-VolumeText 15 "  +4.28%  "

Set-VolumeText 16 "  -0.37%  "

Set-PriceText 17 "0.00001141"
Set-VolumeText 17 "  +1.72%  "

Set-PriceText 18 "93.77"
Set-VolumeText 18 "  +1.67%  "

Set-PriceText 19 "0.06655"
Set-VolumeText 19 "  -0.03%  "

Set-PriceText 20 "18.07"

Set-PriceText 21 "0.9999"
Set-VolumeText 21 "  -0.26%  "

Set-PriceText 22 "6.295"
Set-VolumeText 22 "  +6.59%  "

Set-PriceText 23 "28.622.67"
Set-VolumeText 23 "  +2.35%  "

Set-VolumeText 24 "  +3.13%  "

Set-PriceText 25 "2.280"
Set-VolumeText 25 "  +1.21%  "

Set-PriceText 26 "2.742"
Set-VolumeText 26 "  +15.26%  "

Set-PriceText 27 "2.127.34"
Set-VolumeText 27 "  +5.26%  "

Set-VolumeText 28 "  +3.81%  "

Set-PriceText 29 "159.52"
Set-VolumeText 29 "  +0.23%  "

Set-PriceText 30 "128.93"
Set-VolumeText 30 "  +0.93%  "

Set-PriceText 31 "1.103"
Set-VolumeText 31 "  +6.89%  "

Set-PriceText 32 "0.1074"

Set-PriceText 33 "5.738"
Set-VolumeText 33 "  +3.22%  "

Set-PriceText 34 "3.639"
Set-VolumeText 34 "  +0.07%  "

Set-PriceText 35 "9.891"
Set-VolumeText 35 "  +11.11%  "

Set-PriceText 36 "0.06805"
Set-VolumeText 36 "  +1.61%  "

Set-PriceText 37 "0.02448"
Set-VolumeText 37 "  +5.35%  "

Set-PriceText 38 "1.272"
Set-VolumeText 38 "  +9.84%  "

Set-VolumeText 39 "  +4.68%  "

Set-PriceText 40 "11.75"
Set-VolumeText 40 "  +4.95%  "

Set-PriceText 41 "5.106"
Set-VolumeText 41 "  +3.44%  "

Set-PriceText 42 "0.6442"
Set-VolumeText 42 "  +4.68%  "

Set-PriceText 43 "1.190"
Set-VolumeText 43 "  +4.00%  "

Set-PriceText 45 "13.77"
Set-VolumeText 45 "  +4.59%  "

Set-PriceText 46 "0.6106"
Set-VolumeText 46 "  +3.93%  "

Set-PriceText 47 "1.287"
Set-VolumeText 47 "  -0.45%  "

Set-PriceText 48 "3.677"
Set-VolumeText 48 "  -0.44%  "

Set-VolumeText 49 "  +6.18%  "

Set-PriceText 50 "125.06"
Set-VolumeText 50 "  +1.80%  "

Set-VolumeText 51 "  +3.10%  "
